# Mutual Fund and Registration
# Remove the schedule_type / maximum_amount columns (old U, V) and the
# schedule_config / schedule_verify columns (old Z, AA) from the
# BillPayment data sheet. Deleting whole columns shifts the remaining
# data left so bene_name/bene_query/instrument_type (old W, X, Y) become
# the new U, V, W - matching the target layout exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from right to left so earlier deletions don't shift the
# indices of columns we still need to delete.
$ws.Columns.Item(27).Delete()   # AA: schedule_verify
$ws.Columns.Item(26).Delete()   # Z:  schedule_config
$ws.Columns.Item(22).Delete()   # V:  maximum_amount
$ws.Columns.Item(21).Delete()   # U:  schedule_type
